# Applies the "update description of Tables ; add API for equip" edit
# to the bighw-backend "tables" workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "User" (sheet1)
# ---------------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("User")

# user_id column is now explicitly the student id
$wsUser.Range("C1").Value = "student_id"

# note clarifying what rand_str is used for
$wsUser.Range("I2").Value = "用于验证邮箱与session验证"

# explanatory notes under student_id / user_id, in a small font
$wsUser.Range("C3").Value = "学号，管理员看到"
$wsUser.Range("C4").Value = "一般看到的是默认id(从1开始)"
$wsUser.Range("C3:C4").Font.Size = 6

$wsUser.Columns.Item(3).ColumnWidth = 18.109375

# ---------------------------------------------------------------------
# Sheet "equip" (sheet2)
# ---------------------------------------------------------------------
$wsEquip = $wb.Worksheets.Item("equip")

$wsEquip.Range("A1").Value = "equip_name"
$wsEquip.Range("F2").Value = "可租onsale`n已出租rented`n下架unavailable"

$wsEquip.Range("C3").Value = "去哪领"
$wsEquip.Range("D3").Value = "租期结束时间"
$wsEquip.Range("D3").Font.Size = 8

$wsEquip.Range("A5").Value = "有个默认的id从1开始"

$wsEquip.Columns.Item(1).ColumnWidth = 11.21875
$wsEquip.Columns.Item(2).ColumnWidth = 12.6640625
$wsEquip.Columns.Item(4).ColumnWidth = 10.5546875

$wsEquip.PageSetup.PaperSize = 9
$wsEquip.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Sheet "sale_req" (sheet4) - add lessor_name column
# ---------------------------------------------------------------------
$wsSaleReq = $wb.Worksheets.Item("sale_req")

$wsSaleReq.Range("D1").Value = "lessor_name"
$wsSaleReq.Columns.Item(4).ColumnWidth = 12.6640625

# ---------------------------------------------------------------------
# Sheet "rent_req" (sheet5) - new equip-rental API fields / notes
# ---------------------------------------------------------------------
$wsRentReq = $wb.Worksheets.Item("rent_req")

$wsRentReq.Range("E1").Value = "lessor_name"
$wsRentReq.Range("H6").Value = "申请记录只会被admin删除"
$wsRentReq.Range("H6").Font.Size = 6
$wsRentReq.Range("B8").Value = "默认id从1开始"

$wsRentReq.Columns.Item(5).ColumnWidth = 13.77734375
$wsRentReq.Columns.Item(7).ColumnWidth = 12.44140625

# ---------------------------------------------------------------------
# Selections - restore the cursor positions recorded in each sheet,
# and make "equip" the active tab (it was "rent_req" before).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("auth_req").Range("D2").Select()
$wsSaleReq.Range("B15").Select()
$wsRentReq.Range("B8").Select()
$wb.Worksheets.Item("rent_info").Range("F2").Select()
$wsUser.Range("D6").Select()

$wsEquip.Activate()
$wsEquip.Range("F4").Select()
